# [ST- Updated code For Fixes]
# Apply the same edits that were made to TableCompare/ActualDataFile.xlsx:
#  - Two new rows of "Created" component data were inserted at the top of the
#    "DeveloperTabData" sheet (new testautocomponent_* rows), pushing the
#    previous 3 rows down to rows 3-5.
#  - The newly vacated column F cells on the (now) row 3 & 4 got the same
#    date number format as the rest of column F.
#  - The "DeveloperTabData" tab became the active/selected sheet (previously
#    "Process_SortNode" was selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeveloperTabData")

# Insert two new blank rows at the top; existing rows 1-3 shift down to 3-5.
$ws.Rows("1:2").Insert()

# Populate the new row 2 first, then row 1 -- matches the order the
# original authors entered the data (row 2's new string ends up earlier
# in the shared-strings table than row 1's).
$ws.Range("A2").Value = 341
$ws.Range("B2").Value = "testautocomponent_560547"
$ws.Range("C2").Value = "Shell Script"
$ws.Range("D2").Value = "Suyog Talathi"
$ws.Range("E2").NumberFormat = "m/d/yy h:mm"
$ws.Range("E2").Value = 43244.94027777778
$ws.Range("G2").Value = "Created"

$ws.Range("A1").Value = 345
$ws.Range("B1").Value = "testautocomponent_731106"
$ws.Range("C1").Value = "Shell Script"
$ws.Range("D1").Value = "Suyog Talathi"
$ws.Range("E1").NumberFormat = "m/d/yy h:mm"
$ws.Range("E1").Value = 43245.022847222222
$ws.Range("G1").Value = "Created"

# The old row 3 (A3) and row 4 (A4) gained an (empty) formatted column F
# cell matching the date format used elsewhere in that column.
$ws.Range("F3").NumberFormat = "m/d/yy h:mm"
$ws.Range("F4").NumberFormat = "m/d/yy h:mm"

# Make "DeveloperTabData" the active sheet/tab (was "Process_SortNode").
$ws.Activate()
